# Apply crypto price/volume updates scraped on Sat Jan 13 05:54:31 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.093.85'
$ws.Range("E2").Value = '  -6.62%  '

$ws.Range("D3").Value = '2.553.81'
$ws.Range("E3").Value = '  -1.91%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.48'
$ws.Range("E6").Value = '  -7.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.576'
$ws.Range("E7").Value = '  -3.96%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.553'
$ws.Range("E9").Value = '  -5.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.01'
$ws.Range("E10").Value = '  -8.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0810'
$ws.Range("E11").Value = '  -4.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.75'
$ws.Range("E12").Value = '  -5.38%  '

$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").Value = '2.944.04'
$ws.Range("E14").Value = '  -1.78%  '

$ws.Range("D15").Value = '2.534.75'
$ws.Range("E15").Value = '  -2.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.871'
$ws.Range("E16").Value = '  -5.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.21'
$ws.Range("E17").Value = '  -4.63%  '

$ws.Range("D18").Value = '43.068.46'
$ws.Range("E18").Value = '  -6.90%  '

$ws.Range("E19").Value = '  -1.69%  '

$ws.Range("D20").Value = '0.0₃0980'
$ws.Range("E20").Value = '  -4.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.55'
$ws.Range("E21").Value = '  -2.57%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.21'
$ws.Range("E22").Value = '  -1.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '261.63'
$ws.Range("E23").Value = '  -9.67%  '

$ws.Range("E24").Value = '  -5.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '29.77'
$ws.Range("E25").Value = '  +1.26%  '

$ws.Range("E26").Value = '  -4.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.11'
$ws.Range("E28").Value = '  -7.25%  '

$ws.Range("E29").Value = '  -3.91%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.71'

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.99'
$ws.Range("E31").Value = '  -4.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.01'
$ws.Range("E32").Value = '  -1.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.19'
$ws.Range("E33").Value = '  -0.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.42'
$ws.Range("E34").Value = '  -5.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.73'
$ws.Range("E35").Value = '  -2.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0796'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.116'
$ws.Range("E37").Value = '  -5.49%  '

$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.29'
$ws.Range("E38").Value = '  +14.58%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.119'
$ws.Range("E39").Value = '  -3.56%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.80'
$ws.Range("E40").Value = '  +6.84%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.48'
$ws.Range("E41").Value = '  -3.58%  '

$ws.Range("E42").Value = '  -6.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.85'
$ws.Range("E43").Value = '  -4.78%  '

$ws.Range("D44").Value = '2.084.26'
$ws.Range("E44").Value = '  -1.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '86.01'
$ws.Range("E46").Value = '  -11.78%  '

$ws.Range("E47").Value = '  +3.52%  '

$ws.Range("D48").Value = '2.799.45'
$ws.Range("E48").Value = '  -1.86%  '

$ws.Range("E49").Value = '  -2.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '104.59'
$ws.Range("E50").Value = '  -4.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.69'
$ws.Range("E51").Value = '  -8.51%  '
